# Applies the diff: adds a new "2022-Q1" worksheet (positioned right before
# the "总计" summary sheet) with per-fund holding detail, and inserts a new
# leading row into the "总计" sheet summarizing the 2022-Q1 quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper data: the "2022-Q1" detail rows as they appear in the target
# workbook (fund code, fund name, fund size, stock position, position
# ratio, holding market value (100M yuan), position rank).
# ---------------------------------------------------------------------
$fundRows = @(
    @("005477", "长安鑫禧灵活配置混合A",  "5.70", "94.29", "5.25", "0.2992", 7),
    @("005343", "长安裕盛灵活配置混合A",  "4.66", "94.21", "5.11", "0.2381", 10),
    @("005478", "长安鑫禧灵活配置混合C",  "4.17", "94.29", "5.25", "0.2189", 7),
    @("005344", "长安裕盛灵活配置混合C",  "3.75", "94.21", "5.11", "0.1916", 10),
    @("501030", "汇添富中证环境治理指数（LOF）A", "6.61", "93.20", "2.07", "0.1368", 8),
    @("290014", "泰信现代服务业混合",      "0.73", "81.14", "7.90", "0.0577", 1),
    @("501031", "汇添富中证环境治理指数（LOF）C", "2.74", "93.20", "2.07", "0.0567", 8),
    @("290008", "泰信发展主题混合",        "0.68", "81.03", "7.07", "0.0481", 1),
    @("164908", "交银施罗德中证环境治理指数（LOF）", "2.12", "93.72", "2.10", "0.0445", 8),
    @("000354", "长盛城镇化主题混合",      "0.43", "78.41", "4.52", "0.0194", 10),
    @("010711", "华富国潮优选混合",        "0.25", "94.50", "6.11", "0.0153", 2)
)

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" worksheet, placed immediately before "总计".
# NOTE: worksheet references appear to be positional in this runtime, so
# we must not keep using a cached reference to "总计" after the sheet
# collection has been mutated (e.g. by Add()) - always re-fetch by name.
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$newSheet.Name = "2022-Q1"

# Header row (B1:H1)
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 0; $col -lt $headers.Length; $col++) {
    $cell = $newSheet.Cells.Item(1, $col + 2)
    $cell.Value = $headers[$col]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    foreach ($edge in 7, 8, 9, 10) {
        $border = $cell.Borders.Item($edge)
        $border.LineStyle = 1
        $border.Weight = 2
    }
}

# Data rows (rows 2..12)
for ($i = 0; $i -lt $fundRows.Length; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]

    $idxCell = $newSheet.Cells.Item($r, 1)
    $idxCell.Value = $i
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    foreach ($edge in 7, 8, 9, 10) {
        $border = $idxCell.Borders.Item($edge)
        $border.LineStyle = 1
        $border.Weight = 2
    }

    $codeCell = $newSheet.Cells.Item($r, 2)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $row[0]

    $nameCell = $newSheet.Cells.Item($r, 3)
    $nameCell.NumberFormat = "@"
    $nameCell.Value = $row[1]

    $sizeCell = $newSheet.Cells.Item($r, 4)
    $sizeCell.NumberFormat = "@"
    $sizeCell.Value = $row[2]

    $posCell = $newSheet.Cells.Item($r, 5)
    $posCell.NumberFormat = "@"
    $posCell.Value = $row[3]

    $ratioCell = $newSheet.Cells.Item($r, 6)
    $ratioCell.NumberFormat = "@"
    $ratioCell.Value = $row[4]

    $valueCell = $newSheet.Cells.Item($r, 7)
    $valueCell.NumberFormat = "@"
    $valueCell.Value = $row[5]

    $rankCell = $newSheet.Cells.Item($r, 8)
    $rankCell.Value = $row[6]
}

# ---------------------------------------------------------------------
# 2. Insert a new leading data row into the "总计" sheet for 2022-Q1 and
#    renumber the existing index column. Re-fetch the sheet by name
#    since the worksheet collection changed above.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

$totalIdxCell = $totalSheet.Cells.Item(2, 1)
$totalIdxCell.Value = 0
$totalIdxCell.Font.Bold = $true
$totalIdxCell.HorizontalAlignment = -4108
$totalIdxCell.VerticalAlignment = -4160
foreach ($edge in 7, 8, 9, 10) {
    $border = $totalIdxCell.Borders.Item($edge)
    $border.LineStyle = 1
    $border.Weight = 2
}

$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 11
$totalSheet.Cells.Item(2, 4).Value = 1.33

# Renumber the pre-existing index column (0,1,2 -> 1,2,3) now that a new
# row 0 was inserted above them.
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3
